$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = -11.76
$ws.Range("B7").Value = 4.554199999999995
$ws.Range("A8").Value = -22.40310000000001
$ws.Range("A10").Value = -21.7761
$ws.Range("D10").Value = -7.975199999999998
$ws.Range("A12").Value = -21.53349999999999
$ws.Range("D12").Value = -7.265099999999999
$ws.Range("D13").Value = -9.264399999999997
$ws.Range("D14").Value = -8.008199999999999
$ws.Range("B15").Value = 4.839499999999993
$ws.Range("E16").Value = 16.58130000000001
$ws.Range("A18").Value = -21.59799999999998
$ws.Range("B18").Value = 5.438699999999999
$ws.Range("C18").Value = -10.47519999999999
$ws.Range("E18").Value = 18.26140000000003
$ws.Range("C19").Value = -12.041
$ws.Range("B20").Value = 8.7471
$ws.Range("E22").Value = 16.99280000000001
$ws.Range("E26").Value = 16.34179999999999
$ws.Range("C27").Value = -13.03999999999999
$ws.Range("B29").Value = 4.904300000000002
$ws.Range("D29").Value = -7.484199999999995
$ws.Range("B30").Value = 4.869599999999998
$ws.Range("B31").Value = 4.8678
$ws.Range("C31").Value = -13.5707
$ws.Range("D32").Value = -8.928399999999996
$ws.Range("D35").Value = -7.979
$ws.Range("A37").Value = -19.69439999999999
$ws.Range("C38").Value = -13.20439999999999
$ws.Range("E39").Value = 16.3172
$ws.Range("B40").Value = 9.428999999999995
$ws.Range("C42").Value = -11.8714
$ws.Range("D43").Value = -8.235799999999999
$ws.Range("C44").Value = -13.32349999999999
$ws.Range("E44").Value = 16.40339999999999
$ws.Range("C47").Value = -12.13
$ws.Range("D48").Value = -7.382299999999998
$ws.Range("D49").Value = -8.352000000000006
$ws.Range("B50").Value = 5.461699999999996
$ws.Range("D50").Value = -7.990699999999995
$ws.Range("E51").Value = 17.2791
$ws.Range("E54").Value = 16.50180000000001
$ws.Range("A55").Value = -22.0423
$ws.Range("D56").Value = -8.293399999999998
$ws.Range("E57").Value = 16.66099999999999
$ws.Range("C58").Value = -12.60749999999999
$ws.Range("E63").Value = 18.16120000000001
$ws.Range("C65").Value = -12.3642
$ws.Range("A68").Value = -21.51370000000001
$ws.Range("B68").Value = 4.540700000000002
$ws.Range("D69").Value = -7.190199999999994
$ws.Range("C73").Value = -12.13160000000001
$ws.Range("B76").Value = 6.260699999999998
$ws.Range("A77").Value = -20.11429999999998
$ws.Range("E77").Value = 18.58550000000002
$ws.Range("A78").Value = -19.91679999999998
$ws.Range("A81").Value = -21.7906
$ws.Range("D81").Value = -7.764199999999996
$ws.Range("A82").Value = -21.9918
$ws.Range("E86").Value = 16.75230000000001
$ws.Range("B87").Value = 4.439399999999997
$ws.Range("B88").Value = 4.280599999999999
$ws.Range("C90").Value = -13.1502
$ws.Range("D92").Value = -6.256400000000001
$ws.Range("C94").Value = -10.201
$ws.Range("C95").Value = -12.6326
$ws.Range("B96").Value = 5.604300000000002
$ws.Range("E96").Value = 16.26759999999999
$ws.Range("B98").Value = 5.686499999999999
$ws.Range("E98").Value = 16.1579
$ws.Range("B101").Value = 9.454599999999994
$ws.Range("C101").Value = -12.6007
$ws.Range("B102").Value = 8.756700000000006
